$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.292.87"
$ws.Range("E2").Value = "  -3.48%  "

$ws.Range("D3").Value = "1.932.69"
$ws.Range("E3").Value = "  -3.20%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9994"
$ws.Range("E4").Value = "  +0.05%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "247.74"
$ws.Range("E5").Value = "  -2.74%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.7232"
$ws.Range("E6").Value = "  -10.62%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.9985"
$ws.Range("E7").Value = "  -0.01%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3291"
$ws.Range("E8").Value = "  -6.61%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "26.75"
$ws.Range("E9").Value = "  +4.42%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.06841"
$ws.Range("E10").Value = "  -2.63%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.8096"
$ws.Range("E11").Value = "  -3.84%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07968"
$ws.Range("E12").Value = "  -1.85%  "

$ws.Range("D13").Value = "1.931.54"
$ws.Range("E13").Value = "  -3.21%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.444"
$ws.Range("E14").Value = "  -1.77%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "94.94"
$ws.Range("E15").Value = "  -6.28%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.60"
$ws.Range("E16").Value = "  +4.32%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "263.27"
$ws.Range("E17").Value = "  -3.47%  "

$ws.Range("D18").Value = "30.278.69"
$ws.Range("E18").Value = "  -3.44%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007969"
$ws.Range("E19").Value = "  +0.27%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "5.837"
$ws.Range("E20").Value = "  +0.10%  "

$ws.Range("D21").Value = "2.183.29"
$ws.Range("E21").Value = "  -2.96%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.9988"
$ws.Range("E22").Value = "  -0.01%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.9992"
$ws.Range("E23").Value = "  +0.03%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.923"
$ws.Range("E24").Value = "  -0.77%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.732"
$ws.Range("E25").Value = "  -0.97%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "160.00"
$ws.Range("E26").Value = "  -2.74%  "

$ws.Range("B27").Value = "LidoDAOToken"
$ws.Range("C27").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.340"
$ws.Range("E27").Value = "  +4.79%  "

$ws.Range("B28").Value = "Stellar"
$ws.Range("C28").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.1351"
$ws.Range("E28").Value = "  -10.76%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "19.05"
$ws.Range("E29").Value = "  -5.19%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.368"
$ws.Range("E30").Value = "  +0.70%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.550"
$ws.Range("E31").Value = "  -1.48%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.424"
$ws.Range("E32").Value = "  -3.83%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.224"
$ws.Range("E33").Value = "  -2.96%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.05102"
$ws.Range("E34").Value = "  -2.02%  "

$ws.Range("E35").Value = "  -0.64%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7470"
$ws.Range("E36").Value = "  -1.73%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.725"
$ws.Range("E37").Value = "  -1.74%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01946"
$ws.Range("E38").Value = "  -3.24%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.823"
$ws.Range("E39").Value = "  -3.04%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "80.73"
$ws.Range("E40").Value = "  +2.83%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.607"
$ws.Range("E41").Value = "  -0.67%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.4494"
$ws.Range("E42").Value = "  -5.55%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.020"
$ws.Range("E43").Value = "  -4.89%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.9990"
$ws.Range("E44").Value = "  -0.01%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.8393"
$ws.Range("E45").Value = "  -2.59%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "102.60"
$ws.Range("E46").Value = "  -1.71%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.809"
$ws.Range("E47").Value = "  -2.13%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.349"
$ws.Range("E48").Value = "  -2.44%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "36.38"
$ws.Range("E49").Value = "  -1.54%  "

$ws.Range("B50").Value = "NEARProtocol"
$ws.Range("C50").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.495"
$ws.Range("E50").Value = "  +2.43%  "

$ws.Range("B51").Value = "Decentraland"
$ws.Range("C51").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.4135"
$ws.Range("E51").Value = "  -5.68%  "
